$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Planned" Sprint 3 value (E2): 35 -> 25
$ws.Range("E2").Value = 25

# Update "Actual" Sprint 3 value (E3): 35 -> 25
$ws.Range("E3").Value = 25

# Update "Remaining" Start value (B5): 110 -> 100
# (C5 contains formula =B5-C3, which will recalc to 80 automatically)
$ws.Range("B5").Value = 100

$excel.CalculateFullRebuild()
$wb.Save()
